$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing header cell (AC1) onto the new
# header cells so they match the bold/border/center-top style (s="1")
# used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (W/L/T) repeated for every data row (2-41).
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 83
    $ws.Cells.Item($r, 32).Value = 0
}
